$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5906.7085
$ws.Range("I18").Value = 2752.75
$ws.Range("J18").Value = 6537.5
$ws.Range("K18").Value = 2752.75
$ws.Range("L18").Value = 6537.5
$ws.Range("M18").Value = -2468.75
$ws.Range("N18").Value = -7105.5
$ws.Range("H70").Value = 77779
$ws.Range("J70").Value = 77779
$ws.Range("L70").Value = 233337
$ws.Range("N70").Value = -233877
$ws.Range("H73").Value = 77779
$ws.Range("J73").Value = 77779
$ws.Range("L73").Value = 233337
$ws.Range("N73").Value = -235209
$ws.Range("H74").Value = 7589.533
$ws.Range("I74").Value = 5599.5
$ws.Range("K74").Value = 5599.5
$ws.Range("M74").Value = -4663.5
$ws.Range("H77").Value = 7589.533
$ws.Range("I77").Value = 5599.5
$ws.Range("K77").Value = 27997.5
$ws.Range("M77").Value = -23317.5
$ws.Range("H98").Value = 1507.8889
$ws.Range("I98").Value = 1380.6154
$ws.Range("J98").Value = 1838.8
$ws.Range("K98").Value = 1380.6154
$ws.Range("L98").Value = 1838.8
$ws.Range("M98").Value = 117.3846000000001
$ws.Range("N98").Value = -4834.8
$ws.Range("H107").Value = 515.6667
$ws.Range("I107").Value = 513.7
$ws.Range("K107").Value = 513.7
$ws.Range("M107").Value = 1406.3
$ws.Range("H112").Value = 1245.5883
$ws.Range("J112").Value = 1210.9375
$ws.Range("L112").Value = 3632.8125
$ws.Range("N112").Value = -5848.8125
$ws.Range("H122").Value = 1507.8889
$ws.Range("I122").Value = 1380.6154
$ws.Range("J122").Value = 1838.8
$ws.Range("K122").Value = 4141.8462
$ws.Range("L122").Value = 5516.4
$ws.Range("M122").Value = -1691.8462
$ws.Range("N122").Value = -10416.4
$ws.Range("H129").Value = 3866.625
$ws.Range("I129").Value = 3704.7144
$ws.Range("K129").Value = 11114.1432
$ws.Range("M129").Value = -6114.143199999999
$ws.Range("H132").Value = 4506.476
$ws.Range("I132").Value = 2059.7144
$ws.Range("J132").Value = 9400
$ws.Range("K132").Value = 6179.1432
$ws.Range("L132").Value = 28200
$ws.Range("M132").Value = -3649.1432
$ws.Range("N132").Value = -33260
$ws.Range("H138").Value = 2727.6135
$ws.Range("I138").Value = 1383.5
$ws.Range("J138").Value = 4071.7273
$ws.Range("K138").Value = 4150.5
$ws.Range("L138").Value = 12215.1819
$ws.Range("M138").Value = 989.5
$ws.Range("N138").Value = -22495.1819

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1463.0869
$ws.Range("I32").Value = 816.7595
$ws.Range("J32").Value = 5390.769
$ws.Range("K32").Value = 816.7595
$ws.Range("L32").Value = 5390.769
$ws.Range("M32").Value = -529.7595
$ws.Range("N32").Value = -5964.769
$ws.Range("H61").Value = 1770.0256
$ws.Range("I61").Value = 1522.9062
$ws.Range("K61").Value = 1522.9062
$ws.Range("M61").Value = -1310.9062
$ws.Range("H110").Value = 2287.6667
$ws.Range("I110").Value = 2145.2
$ws.Range("K110").Value = 2145.2
$ws.Range("M110").Value = -100.1999999999998
$ws.Range("H132").Value = 2895.25
$ws.Range("I132").Value = 2690.2354
$ws.Range("K132").Value = 8070.706200000001
$ws.Range("M132").Value = -5540.706200000001
$ws.Range("H136").Value = 1770.0256
$ws.Range("I136").Value = 1522.9062
$ws.Range("K136").Value = 4568.7186
$ws.Range("M136").Value = -2018.7186
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 878.1
$ws.Range("I20").Value = 792.7895
$ws.Range("J20").Value = 1025.4546
$ws.Range("K20").Value = 792.7895
$ws.Range("L20").Value = 1025.4546
$ws.Range("M20").Value = -545.7895
$ws.Range("N20").Value = -1519.4546
$ws.Range("H86").Value = 1642
$ws.Range("J86").Value = 1713
$ws.Range("L86").Value = 1713
$ws.Range("N86").Value = -3959
$ws.Range("H89").Value = 1642
$ws.Range("J89").Value = 1713
$ws.Range("L89").Value = 8565
$ws.Range("N89").Value = -19797
$ws.Range("H94").Value = 1128.2609
$ws.Range("I94").Value = 1093.5883
$ws.Range("J94").Value = 1226.5
$ws.Range("K94").Value = 1093.5883
$ws.Range("L94").Value = 1226.5
$ws.Range("M94").Value = -642.5882999999999
$ws.Range("N94").Value = -2128.5
$ws.Range("H107").Value = 9826.857
$ws.Range("I107").Value = 11049.667
$ws.Range("K107").Value = 11049.667
$ws.Range("M107").Value = -9129.666999999999
$ws.Range("H134").Value = 2062.6196
$ws.Range("I134").Value = 1236.5151
$ws.Range("J134").Value = 4159.654
$ws.Range("K134").Value = 3709.5453
$ws.Range("L134").Value = 12478.962
$ws.Range("M134").Value = -1174.5453
$ws.Range("N134").Value = -17548.962

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 948.8
$ws.Range("I22").Value = 871.2222
$ws.Range("K22").Value = 871.2222
$ws.Range("M22").Value = -521.2222
$ws.Range("H132").Value = 1595.6111
$ws.Range("I132").Value = 1595.6111
$ws.Range("K132").Value = 4786.8333
$ws.Range("M132").Value = -2256.8333
$ws.Range("H134").Value = 2055.4
$ws.Range("I134").Value = 1333.375
$ws.Range("J134").Value = 3832.6924
$ws.Range("K134").Value = 4000.125
$ws.Range("L134").Value = 11498.0772
$ws.Range("M134").Value = -1465.125
$ws.Range("N134").Value = -16568.0772

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1187.4445
$ws.Range("I29").Value = 310.8
$ws.Range("K29").Value = 932.4000000000001
$ws.Range("M29").Value = -655.4000000000001
$ws.Range("H37").Value = 97662.164
$ws.Range("J37").Value = 97662.164
$ws.Range("L37").Value = 292986.492
$ws.Range("N37").Value = -293210.492
$ws.Range("H117").Value = 5478.769
$ws.Range("J117").Value = 6524.8887
$ws.Range("L117").Value = 19574.6661
$ws.Range("N117").Value = -26458.6661
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H131").Value = 1199.7391
$ws.Range("J131").Value = 1240.5476
$ws.Range("L131").Value = 3721.642800000001
$ws.Range("N131").Value = -13801.6428
$ws.Range("H137").Value = 4207.5
$ws.Range("I137").Value = 2560.8333
$ws.Range("J137").Value = 4756.3887
$ws.Range("K137").Value = 7682.499899999999
$ws.Range("L137").Value = 14269.1661
$ws.Range("M137").Value = -2582.499899999999
$ws.Range("N137").Value = -24469.1661
$ws.Range("H140").Value = 3186
$ws.Range("I140").Value = 2634.182
$ws.Range("J140").Value = 4400
$ws.Range("K140").Value = 7902.545999999999
$ws.Range("L140").Value = 13200
$ws.Range("M140").Value = -2722.545999999999
$ws.Range("N140").Value = -23560

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 49951
$ws.Range("J93").Value = 49951
$ws.Range("L93").Value = 49951
$ws.Range("N93").Value = -53695
$ws.Range("H102").Value = 2570.6072
$ws.Range("I102").Value = 2554.7036
$ws.Range("K102").Value = 2554.7036
$ws.Range("M102").Value = -932.7035999999998
$ws.Range("H132").Value = 27037414
$ws.Range("I132").Value = 34490252
$ws.Range("J132").Value = 20874.125
$ws.Range("K132").Value = 103470756
$ws.Range("L132").Value = 62622.375
$ws.Range("M132").Value = -103468226
$ws.Range("N132").Value = -67682.375
$ws.Range("H135").Value = 95000
$ws.Range("J135").Value = 95000
$ws.Range("L135").Value = 95000
$ws.Range("N135").Value = -105140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4334
$ws.Range("I61").Value = 1199.75
$ws.Range("J61").Value = 10602.5
$ws.Range("K61").Value = 1199.75
$ws.Range("L61").Value = 10602.5
$ws.Range("M61").Value = -997.75
$ws.Range("N61").Value = -11006.5
$ws.Range("H93").Value = 2372.8
$ws.Range("I93").Value = 2299.3333
$ws.Range("J93").Value = 2666.6667
$ws.Range("K93").Value = 2299.3333
$ws.Range("L93").Value = 2666.6667
$ws.Range("M93").Value = -1051.3333
$ws.Range("N93").Value = -5162.6667
$ws.Range("H113").Value = 4334
$ws.Range("I113").Value = 1199.75
$ws.Range("J113").Value = 10602.5
$ws.Range("K113").Value = 1199.75
$ws.Range("L113").Value = 10602.5
$ws.Range("M113").Value = 970.25
$ws.Range("N113").Value = -14942.5
$ws.Range("H122").Value = 4451.8335
$ws.Range("I122").Value = 3668.6538
$ws.Range("K122").Value = 11005.9614
$ws.Range("M122").Value = -8555.9614

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2103.0264
$ws.Range("I107").Value = 1193.52
$ws.Range("K107").Value = 3580.56
$ws.Range("M107").Value = -1660.56
$ws.Range("H132").Value = 2597.6538
$ws.Range("I132").Value = 1749
$ws.Range("K132").Value = 5247
$ws.Range("M132").Value = -2717
